$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Daily refresh of the COVID-19 "paises" tracker ---

# 1) Update the "last updated" timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Agosto de 2020 a las 01:09"

# 2) A handful of countries changed rank (sorted by "Casos totales" desc),
#    so their name swaps with the neighboring row.
$countrySwaps = @(
    @{ Row = 12;  Name = "Colombia" },
    @{ Row = 13;  Name = "España" },
    @{ Row = 50;  Name = "Nigeria" },
    @{ Row = 51;  Name = "Honduras" },
    @{ Row = 74;  Name = "Chequia" },
    @{ Row = 75;  Name = "Camerun" },
    @{ Row = 202; Name = "Timor Oriental" },
    @{ Row = 203; Name = "Santa Lucia" }
)
foreach ($swap in $countrySwaps) {
    $ws.Cells.Item($swap.Row, 1).Value = $swap.Name
}

# 3) Updated per-country statistics: Casos totales(B), Nuevos casos(C),
#    Casos activos(D), Recuperados(E), Casos criticos(F), Muertes hoy(G), Muertes(H)
$statUpdates = @(
    @{ Row = 4; B = 5027191; C = 53623; D = 2569431; E = 2295071; G = 1088; H = 162689 },
    @{ Row = 5; B = 2917562; C = 54801; D = 2047660; E = 771258; G = 1226; H = 98644 },
    @{ Row = 12; B = 357710; C = 11996; D = 192355; E = 153416; G = 315; H = 11939 },
    @{ Row = 13; B = 354530; C = 1683; D = 0; E = 0; G = 1; H = 28500 },
    @{ Row = 27; B = 118514; C = 327; D = 103077; E = 6471; G = 4; H = 8966 },
    @{ Row = 31; B = 90537; C = 1671; D = 71318; E = 13342; G = 30; H = 5877 },
    @{ Row = 50; B = 45244; C = 354; D = 32430; E = 11884; G = 3; H = 930 },
    @{ Row = 51; B = 45098; C = 799; D = 6116; E = 37559; G = 23; H = 1423 },
    @{ Row = 52; B = 42889; C = 375; D = 39945; E = 2788 },
    @{ Row = 53; B = 42263; C = 1134; D = 28877; E = 12360; G = 4; H = 1026 },
    @{ Row = 74; B = 17731; C = 202; D = 12320; E = 5021; G = 2; H = 390 },
    @{ Row = 75; B = 17718; C = 0; D = 15320; E = 2007; G = 0; H = 391 },
    @{ Row = 81; B = 13014; C = 297; D = 7374; E = 5205; G = 11; H = 435 },
    @{ Row = 86; B = 9468; C = 59; E = 355 },
    @{ Row = 130; B = 2111; C = 7; D = 1258; E = 848 },
    @{ Row = 136; B = 1768; C = 5; D = 898; E = 362 },
    @{ Row = 140; B = 1318; C = 9; D = 1079; E = 202 },
    @{ Row = 153; B = 942; C = 3; D = 838; E = 28; G = 1; H = 76 },
    @{ Row = 179; B = 210; C = 11; E = 67 }
)
foreach ($u in $statUpdates) {
    if ($u.ContainsKey("B")) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey("C")) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey("D")) { $ws.Cells.Item($u.Row, 4).Value = $u.D }
    if ($u.ContainsKey("E")) { $ws.Cells.Item($u.Row, 5).Value = $u.E }
    if ($u.ContainsKey("F")) { $ws.Cells.Item($u.Row, 6).Value = $u.F }
    if ($u.ContainsKey("G")) { $ws.Cells.Item($u.Row, 7).Value = $u.G }
    if ($u.ContainsKey("H")) { $ws.Cells.Item($u.Row, 8).Value = $u.H }
}
